$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing Row 2 (Beta) values for columns F:N
$ws.Range("F2").Value = 11.87521617810516
$ws.Range("G2").Value = 11.64411788020399
$ws.Range("H2").Value = 12.10248894967641
$ws.Range("I2").Value = 1.963881362547014
$ws.Range("J2").Value = 1.942777076614074
$ws.Range("K2").Value = 1.985340881616313
$ws.Range("L2").Value = 0.1527875343608597
$ws.Range("M2").Value = 0.1511420402856856
$ws.Range("N2").Value = 0.1544346540239098

# Update existing Row 3 (Gamma) values for columns F:N
$ws.Range("F3").Value = 0.001952783590897234
$ws.Range("G3").Value = 0.00120252990512603
$ws.Range("H3").Value = 0.002847997915831956
$ws.Range("I3").Value = 0.001808441993160187
$ws.Range("J3").Value = 0.001105933764838022
$ws.Range("K3").Value = 0.002644212699523138
$ws.Range("L3").Value = 0.002036737347825657
$ws.Range("M3").Value = 0.00127183433226553
$ws.Range("N3").Value = 0.002946034071592752

# Add new Row 4 (Beta + Gamma)
$ws.Range("A2").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Beta + Gamma"
$ws.Range("C4").Value = 12.00687180793019
$ws.Range("D4").Value = 1.974477778970852
$ws.Range("E4").Value = 0.1537386519519979
$ws.Range("F4").Value = 11.87716896169606
$ws.Range("G4").Value = 11.64532041010911
$ws.Range("H4").Value = 12.10533694759224
$ws.Range("I4").Value = 1.965689804540174
$ws.Range("J4").Value = 1.943883010378911
$ws.Range("K4").Value = 1.987985094315837
$ws.Range("L4").Value = 0.1548242717086853
$ws.Range("M4").Value = 0.1524138746179511
$ws.Range("N4").Value = 0.1573806880955025
